# Apply the answers-of-within100 update (commit 9703e77):
# - bump the header date from 2023-03-27 Monday to 2023-03-28 Tuesday
# - regenerate each of the 100 arithmetic "answer" cells in the table
#
# Every "old" value below is unique in the document, so a plain
# Find/Replace (first match, case-sensitive, whole document) lands on
# exactly the intended run each time. Replacements are issued in the
# same order the cells appear in the document so that the one value
# that is reused later as a *new* answer ("16-13=3") is only ever
# matched against its original occurrence before that occurrence is
# itself replaced.

$d = $word.ActiveDocument

# Header date
$null = $d.Content.Find.Execute("2023-03-27 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-28 Tuesday", 2)

# Table answer cells, in document order
$null = $d.Content.Find.Execute("43+43=86", $true, $false, $false, $false, $false, $true, 1, $false, "81+0=81", 2)
$null = $d.Content.Find.Execute("22+53=75", $true, $false, $false, $false, $false, $true, 1, $false, "64+32=96", 2)
$null = $d.Content.Find.Execute("39-21=18", $true, $false, $false, $false, $false, $true, 1, $false, "30-2=28", 2)
$null = $d.Content.Find.Execute("75-54=21", $true, $false, $false, $false, $false, $true, 1, $false, "4+33=37", 2)
$null = $d.Content.Find.Execute("31-21=10", $true, $false, $false, $false, $false, $true, 1, $false, "31-4=27", 2)
$null = $d.Content.Find.Execute("13+36=49", $true, $false, $false, $false, $false, $true, 1, $false, "93-36=57", 2)
$null = $d.Content.Find.Execute("42-36=6", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=68", 2)
$null = $d.Content.Find.Execute("33+14=47", $true, $false, $false, $false, $false, $true, 1, $false, "68-43=25", 2)
$null = $d.Content.Find.Execute("67-58=9", $true, $false, $false, $false, $false, $true, 1, $false, "7+24=31", 2)
$null = $d.Content.Find.Execute("43+33=76", $true, $false, $false, $false, $false, $true, 1, $false, "86-14=72", 2)
$null = $d.Content.Find.Execute("10+71=81", $true, $false, $false, $false, $false, $true, 1, $false, "30-2=28", 2)
$null = $d.Content.Find.Execute("77-15=62", $true, $false, $false, $false, $false, $true, 1, $false, "5+7=12", 2)
$null = $d.Content.Find.Execute("25+13=38", $true, $false, $false, $false, $false, $true, 1, $false, "39+34=73", 2)
$null = $d.Content.Find.Execute("89-39=50", $true, $false, $false, $false, $false, $true, 1, $false, "42+49=91", 2)
$null = $d.Content.Find.Execute("52+7=59", $true, $false, $false, $false, $false, $true, 1, $false, "22+66=88", 2)
$null = $d.Content.Find.Execute("6+88=94", $true, $false, $false, $false, $false, $true, 1, $false, "48+30=78", 2)
$null = $d.Content.Find.Execute("61-43=18", $true, $false, $false, $false, $false, $true, 1, $false, "46-34=12", 2)
$null = $d.Content.Find.Execute("10+70=80", $true, $false, $false, $false, $false, $true, 1, $false, "51-32=19", 2)
$null = $d.Content.Find.Execute("24+56=80", $true, $false, $false, $false, $false, $true, 1, $false, "31+64=95", 2)
$null = $d.Content.Find.Execute("25-2=23", $true, $false, $false, $false, $false, $true, 1, $false, "4+32=36", 2)
$null = $d.Content.Find.Execute("52-6=46", $true, $false, $false, $false, $false, $true, 1, $false, "11+75=86", 2)
$null = $d.Content.Find.Execute("21-8=13", $true, $false, $false, $false, $false, $true, 1, $false, "77-74=3", 2)
$null = $d.Content.Find.Execute("21+74=95", $true, $false, $false, $false, $false, $true, 1, $false, "56-34=22", 2)
$null = $d.Content.Find.Execute("43+36=79", $true, $false, $false, $false, $false, $true, 1, $false, "99-68=31", 2)
$null = $d.Content.Find.Execute("16-13=3", $true, $false, $false, $false, $false, $true, 1, $false, "46-15=31", 2)
$null = $d.Content.Find.Execute("35+64=99", $true, $false, $false, $false, $false, $true, 1, $false, "34+64=98", 2)
$null = $d.Content.Find.Execute("76-54=22", $true, $false, $false, $false, $false, $true, 1, $false, "16-13=3", 2)
$null = $d.Content.Find.Execute("25+35=60", $true, $false, $false, $false, $false, $true, 1, $false, "17+68=85", 2)
$null = $d.Content.Find.Execute("25+72=97", $true, $false, $false, $false, $false, $true, 1, $false, "90-67=23", 2)
$null = $d.Content.Find.Execute("5+65=70", $true, $false, $false, $false, $false, $true, 1, $false, "85-63=22", 2)
$null = $d.Content.Find.Execute("45+21=66", $true, $false, $false, $false, $false, $true, 1, $false, "72-11=61", 2)
$null = $d.Content.Find.Execute("0+31=31", $true, $false, $false, $false, $false, $true, 1, $false, "13+3=16", 2)
$null = $d.Content.Find.Execute("80-49=31", $true, $false, $false, $false, $false, $true, 1, $false, "99-0=99", 2)
$null = $d.Content.Find.Execute("65-41=24", $true, $false, $false, $false, $false, $true, 1, $false, "76+23=99", 2)
$null = $d.Content.Find.Execute("7+1=8", $true, $false, $false, $false, $false, $true, 1, $false, "9+57=66", 2)
$null = $d.Content.Find.Execute("87-85=2", $true, $false, $false, $false, $false, $true, 1, $false, "29+64=93", 2)
$null = $d.Content.Find.Execute("72-28=44", $true, $false, $false, $false, $false, $true, 1, $false, "65-58=7", 2)
$null = $d.Content.Find.Execute("38+16=54", $true, $false, $false, $false, $false, $true, 1, $false, "46+52=98", 2)
$null = $d.Content.Find.Execute("39+29=68", $true, $false, $false, $false, $false, $true, 1, $false, "9+90=99", 2)
$null = $d.Content.Find.Execute("81-69=12", $true, $false, $false, $false, $false, $true, 1, $false, "7+5=12", 2)
$null = $d.Content.Find.Execute("66+31=97", $true, $false, $false, $false, $false, $true, 1, $false, "19+32=51", 2)
$null = $d.Content.Find.Execute("80-61=19", $true, $false, $false, $false, $false, $true, 1, $false, "82-72=10", 2)
$null = $d.Content.Find.Execute("21+37=58", $true, $false, $false, $false, $false, $true, 1, $false, "66+25=91", 2)
$null = $d.Content.Find.Execute("57-25=32", $true, $false, $false, $false, $false, $true, 1, $false, "60+16=76", 2)
$null = $d.Content.Find.Execute("23+74=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-37=52", 2)
$null = $d.Content.Find.Execute("88-9=79", $true, $false, $false, $false, $false, $true, 1, $false, "14+42=56", 2)
$null = $d.Content.Find.Execute("48+22=70", $true, $false, $false, $false, $false, $true, 1, $false, "60+32=92", 2)
$null = $d.Content.Find.Execute("20+68=88", $true, $false, $false, $false, $false, $true, 1, $false, "52-2=50", 2)
$null = $d.Content.Find.Execute("74-7=67", $true, $false, $false, $false, $false, $true, 1, $false, "86-66=20", 2)
$null = $d.Content.Find.Execute("64-56=8", $true, $false, $false, $false, $false, $true, 1, $false, "56-2=54", 2)
$null = $d.Content.Find.Execute("40-30=10", $true, $false, $false, $false, $false, $true, 1, $false, "47+45=92", 2)
$null = $d.Content.Find.Execute("62-58=4", $true, $false, $false, $false, $false, $true, 1, $false, "97-47=50", 2)
$null = $d.Content.Find.Execute("69-61=8", $true, $false, $false, $false, $false, $true, 1, $false, "50+23=73", 2)
$null = $d.Content.Find.Execute("6+2=8", $true, $false, $false, $false, $false, $true, 1, $false, "35-22=13", 2)
$null = $d.Content.Find.Execute("90-74=16", $true, $false, $false, $false, $false, $true, 1, $false, "73+24=97", 2)
$null = $d.Content.Find.Execute("89+7=96", $true, $false, $false, $false, $false, $true, 1, $false, "53-32=21", 2)
$null = $d.Content.Find.Execute("2+20=22", $true, $false, $false, $false, $false, $true, 1, $false, "46-34=12", 2)
$null = $d.Content.Find.Execute("77+2=79", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=91", 2)
$null = $d.Content.Find.Execute("16-6=10", $true, $false, $false, $false, $false, $true, 1, $false, "76+17=93", 2)
$null = $d.Content.Find.Execute("94-50=44", $true, $false, $false, $false, $false, $true, 1, $false, "59+9=68", 2)
$null = $d.Content.Find.Execute("30+10=40", $true, $false, $false, $false, $false, $true, 1, $false, "77-49=28", 2)
$null = $d.Content.Find.Execute("81-38=43", $true, $false, $false, $false, $false, $true, 1, $false, "60+35=95", 2)
$null = $d.Content.Find.Execute("12+38=50", $true, $false, $false, $false, $false, $true, 1, $false, "13-3=10", 2)
$null = $d.Content.Find.Execute("2+97=99", $true, $false, $false, $false, $false, $true, 1, $false, "22+28=50", 2)
$null = $d.Content.Find.Execute("25-22=3", $true, $false, $false, $false, $false, $true, 1, $false, "70-52=18", 2)
$null = $d.Content.Find.Execute("48+16=64", $true, $false, $false, $false, $false, $true, 1, $false, "13-3=10", 2)
$null = $d.Content.Find.Execute("3+76=79", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=76", 2)
$null = $d.Content.Find.Execute("9-5=4", $true, $false, $false, $false, $false, $true, 1, $false, "54+37=91", 2)
$null = $d.Content.Find.Execute("22+41=63", $true, $false, $false, $false, $false, $true, 1, $false, "86-83=3", 2)
$null = $d.Content.Find.Execute("39+2=41", $true, $false, $false, $false, $false, $true, 1, $false, "87-83=4", 2)
$null = $d.Content.Find.Execute("76+16=92", $true, $false, $false, $false, $false, $true, 1, $false, "69+11=80", 2)
$null = $d.Content.Find.Execute("41-12=29", $true, $false, $false, $false, $false, $true, 1, $false, "43-19=24", 2)
$null = $d.Content.Find.Execute("92+6=98", $true, $false, $false, $false, $false, $true, 1, $false, "21-14=7", 2)
$null = $d.Content.Find.Execute("83-73=10", $true, $false, $false, $false, $false, $true, 1, $false, "61+17=78", 2)
$null = $d.Content.Find.Execute("96-67=29", $true, $false, $false, $false, $false, $true, 1, $false, "95-44=51", 2)
$null = $d.Content.Find.Execute("62-30=32", $true, $false, $false, $false, $false, $true, 1, $false, "1+3=4", 2)
$null = $d.Content.Find.Execute("74+12=86", $true, $false, $false, $false, $false, $true, 1, $false, "13-2=11", 2)
$null = $d.Content.Find.Execute("96-60=36", $true, $false, $false, $false, $false, $true, 1, $false, "92-81=11", 2)
$null = $d.Content.Find.Execute("31-16=15", $true, $false, $false, $false, $false, $true, 1, $false, "95-18=77", 2)
$null = $d.Content.Find.Execute("19+45=64", $true, $false, $false, $false, $false, $true, 1, $false, "28-4=24", 2)
$null = $d.Content.Find.Execute("77+7=84", $true, $false, $false, $false, $false, $true, 1, $false, "55+39=94", 2)
$null = $d.Content.Find.Execute("2+65=67", $true, $false, $false, $false, $false, $true, 1, $false, "51+37=88", 2)
$null = $d.Content.Find.Execute("53-17=36", $true, $false, $false, $false, $false, $true, 1, $false, "54-41=13", 2)
$null = $d.Content.Find.Execute("14+80=94", $true, $false, $false, $false, $false, $true, 1, $false, "61-25=36", 2)
$null = $d.Content.Find.Execute("76-6=70", $true, $false, $false, $false, $false, $true, 1, $false, "15+18=33", 2)
$null = $d.Content.Find.Execute("40+10=50", $true, $false, $false, $false, $false, $true, 1, $false, "21+32=53", 2)
$null = $d.Content.Find.Execute("0+22=22", $true, $false, $false, $false, $false, $true, 1, $false, "56-37=19", 2)
$null = $d.Content.Find.Execute("41+35=76", $true, $false, $false, $false, $false, $true, 1, $false, "52+0=52", 2)
$null = $d.Content.Find.Execute("12+22=34", $true, $false, $false, $false, $false, $true, 1, $false, "95-53=42", 2)
$null = $d.Content.Find.Execute("84-35=49", $true, $false, $false, $false, $false, $true, 1, $false, "2+27=29", 2)
$null = $d.Content.Find.Execute("84-2=82", $true, $false, $false, $false, $false, $true, 1, $false, "70+14=84", 2)
$null = $d.Content.Find.Execute("63-49=14", $true, $false, $false, $false, $false, $true, 1, $false, "36+53=89", 2)
$null = $d.Content.Find.Execute("98-7=91", $true, $false, $false, $false, $false, $true, 1, $false, "14+25=39", 2)
$null = $d.Content.Find.Execute("81-13=68", $true, $false, $false, $false, $false, $true, 1, $false, "59-32=27", 2)
$null = $d.Content.Find.Execute("20+43=63", $true, $false, $false, $false, $false, $true, 1, $false, "40-27=13", 2)
$null = $d.Content.Find.Execute("5+21=26", $true, $false, $false, $false, $false, $true, 1, $false, "83-29=54", 2)
$null = $d.Content.Find.Execute("17+78=95", $true, $false, $false, $false, $false, $true, 1, $false, "32+14=46", 2)
$null = $d.Content.Find.Execute("70+22=92", $true, $false, $false, $false, $false, $true, 1, $false, "60-40=20", 2)
$null = $d.Content.Find.Execute("39+20=59", $true, $false, $false, $false, $false, $true, 1, $false, "63+21=84", 2)
$null = $d.Content.Find.Execute("19+29=48", $true, $false, $false, $false, $false, $true, 1, $false, "99-97=2", 2)
